$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H15").Value = 3410.5193
$ws_ALC.Range("I15").Value = 3410.5193
$ws_ALC.Range("K15").Value = 10231.5579
$ws_ALC.Range("M15").Value = -10062.5579

$ws_ALC.Range("H51").Value = 5800.3335
$ws_ALC.Range("I51").Value = 4978.5835
$ws_ALC.Range("J51").Value = 7443.8335
$ws_ALC.Range("K51").Value = 4978.5835
$ws_ALC.Range("L51").Value = 7443.8335
$ws_ALC.Range("M51").Value = -4494.5835
$ws_ALC.Range("N51").Value = -8411.833500000001

$ws_ALC.Range("H62").Value = 3057.7
$ws_ALC.Range("I62").Value = 3008.5557
$ws_ALC.Range("K62").Value = 3008.5557
$ws_ALC.Range("M62").Value = -2384.5557

$ws_ALC.Range("H65").Value = 3057.7
$ws_ALC.Range("I65").Value = 3008.5557
$ws_ALC.Range("K65").Value = 15042.7785
$ws_ALC.Range("M65").Value = -11922.7785

$ws_ALC.Range("H74").Value = 15350.292
$ws_ALC.Range("I74").Value = 16094.833
$ws_ALC.Range("K74").Value = 16094.833
$ws_ALC.Range("M74").Value = -15158.833

$ws_ALC.Range("H76").Value = 4331.3335
$ws_ALC.Range("I76").Value = 4000
$ws_ALC.Range("K76").Value = 4000
$ws_ALC.Range("M76").Value = -3685

$ws_ALC.Range("H77").Value = 15350.292
$ws_ALC.Range("I77").Value = 16094.833
$ws_ALC.Range("K77").Value = 80474.16500000001
$ws_ALC.Range("M77").Value = -75794.16500000001

$ws_ALC.Range("H79").Value = 4331.3335
$ws_ALC.Range("I79").Value = 4000
$ws_ALC.Range("K79").Value = 4000
$ws_ALC.Range("M79").Value = -2908

$ws_ALC.Range("H112").Value = 2026.8966
$ws_ALC.Range("J112").Value = 2062.2222
$ws_ALC.Range("L112").Value = 6186.6666
$ws_ALC.Range("N112").Value = -8402.6666

$ws_ALC.Range("H115").Value = 2941.6
$ws_ALC.Range("I115").Value = 402.66666
$ws_ALC.Range("K115").Value = 1207.99998
$ws_ALC.Range("M115").Value = 359.0000199999999

$ws_ALC.Range("H116").Value = 5499.75
$ws_ALC.Range("I116").Value = 3999.5
$ws_ALC.Range("J116").Value = 7000
$ws_ALC.Range("K116").Value = 3999.5
$ws_ALC.Range("L116").Value = 7000
$ws_ALC.Range("M116").Value = -557.5
$ws_ALC.Range("N116").Value = -13884

$ws_ALC.Range("H132").Value = 6381.933
$ws_ALC.Range("I132").Value = 6878.48
$ws_ALC.Range("K132").Value = 20635.44
$ws_ALC.Range("M132").Value = -18105.44

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H25").Value = 11000
$ws_ARM.Range("I25").Value = 11000
$ws_ARM.Range("K25").Value = 11000
$ws_ARM.Range("M25").Value = -10598

$ws_ARM.Range("H32").Value = 3977.9
$ws_ARM.Range("I32").Value = 3339.0977
$ws_ARM.Range("J32").Value = 6888
$ws_ARM.Range("K32").Value = 3339.0977
$ws_ARM.Range("L32").Value = 6888
$ws_ARM.Range("M32").Value = -3052.0977
$ws_ARM.Range("N32").Value = -7462

$ws_ARM.Range("H45").Value = 23607.578
$ws_ARM.Range("I45").Value = 29216.2
$ws_ARM.Range("K45").Value = 29216.2
$ws_ARM.Range("M45").Value = -28839.2

$ws_ARM.Range("H132").Value = 1772.9166
$ws_ARM.Range("I132").Value = 1160.8667
$ws_ARM.Range("J132").Value = 4833.1665
$ws_ARM.Range("K132").Value = 3482.6001
$ws_ARM.Range("L132").Value = 14499.4995
$ws_ARM.Range("M132").Value = -952.6001000000001
$ws_ARM.Range("N132").Value = -19559.4995

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H31").Value = 5999
$ws_BSM.Range("I31").Value = 5999
$ws_BSM.Range("K31").Value = 5999
$ws_BSM.Range("M31").Value = -5747

$ws_BSM.Range("H86").Value = 2066.611
$ws_BSM.Range("I86").Value = 1809.5385
$ws_BSM.Range("K86").Value = 1809.5385
$ws_BSM.Range("M86").Value = -686.5385000000001

$ws_BSM.Range("H89").Value = 2066.611
$ws_BSM.Range("I89").Value = 1809.5385
$ws_BSM.Range("K89").Value = 9047.692500000001
$ws_BSM.Range("M89").Value = -3431.692500000001

$ws_BSM.Range("H94").Value = 90915300
$ws_BSM.Range("I94").Value = 142864800
$ws_BSM.Range("K94").Value = 142864800
$ws_BSM.Range("M94").Value = -142864349

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H23").Value = 5500000
$ws_CRP.Range("I23").Value = 5500000
$ws_CRP.Range("J23").Value = 0
$ws_CRP.Range("K23").Value = 5500000
$ws_CRP.Range("L23").Value = 0
$ws_CRP.Range("M23").ClearContents()
$ws_CRP.Range("N23").Value = -5499760

$ws_CRP.Range("H27").Value = 5500000
$ws_CRP.Range("I27").Value = 5500000
$ws_CRP.Range("J27").Value = 0
$ws_CRP.Range("K27").Value = 5500000
$ws_CRP.Range("L27").Value = 0
$ws_CRP.Range("M27").ClearContents()
$ws_CRP.Range("N27").Value = -5499808

$ws_CRP.Range("H58").Value = 2241.524
$ws_CRP.Range("I58").Value = 1191.091
$ws_CRP.Range("K58").Value = 1191.091
$ws_CRP.Range("M58").Value = -988.0909999999999

$ws_CRP.Range("H122").Value = 4479.4165
$ws_CRP.Range("I122").Value = 3596.2307
$ws_CRP.Range("J122").Value = 5523.1816
$ws_CRP.Range("K122").Value = 10788.6921
$ws_CRP.Range("L122").Value = 16569.5448
$ws_CRP.Range("M122").Value = -8338.6921
$ws_CRP.Range("N122").Value = -21469.5448

$ws_CRP.Range("H136").Value = 2241.524
$ws_CRP.Range("I136").Value = 1191.091
$ws_CRP.Range("K136").Value = 3573.273
$ws_CRP.Range("M136").Value = -1023.273

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H137").Value = 83077.664
$ws_GSM.Range("J137").Value = 81693.2
$ws_GSM.Range("L137").Value = 81693.2
$ws_GSM.Range("N137").Value = -91893.2

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 642.9167
$ws_LTW.Range("I22").Value = 545.1429000000001
$ws_LTW.Range("K22").Value = 545.1429000000001
$ws_LTW.Range("M22").Value = -250.1429000000001

$ws_LTW.Range("H27").Value = 642.9167
$ws_LTW.Range("I27").Value = 545.1429000000001
$ws_LTW.Range("K27").Value = 545.1429000000001
$ws_LTW.Range("M27").Value = -438.1429000000001

$ws_LTW.Range("H40").Value = 4692.6733
$ws_LTW.Range("I40").Value = 4568.72
$ws_LTW.Range("K40").Value = 4568.72
$ws_LTW.Range("M40").Value = -4432.72

$ws_LTW.Range("H82").Value = 1167
$ws_LTW.Range("I82").Value = 1449.8
$ws_LTW.Range("J82").Value = 695.6667
$ws_LTW.Range("K82").Value = 1449.8
$ws_LTW.Range("L82").Value = 695.6667
$ws_LTW.Range("M82").Value = -1088.8
$ws_LTW.Range("N82").Value = -1417.6667

$ws_LTW.Range("H85").Value = 1167
$ws_LTW.Range("I85").Value = 1449.8
$ws_LTW.Range("J85").Value = 695.6667
$ws_LTW.Range("K85").Value = 1449.8
$ws_LTW.Range("L85").Value = 695.6667
$ws_LTW.Range("M85").Value = -201.8
$ws_LTW.Range("N85").Value = -3191.6667

$ws_LTW.Range("H100").Value = 1704
$ws_LTW.Range("I100").Value = 1399.5
$ws_LTW.Range("K100").Value = 1399.5
$ws_LTW.Range("M100").Value = -858.5

$ws_LTW.Range("H134").Value = 103981.336
$ws_LTW.Range("J134").Value = 103981.336
$ws_LTW.Range("L134").Value = 103981.336
$ws_LTW.Range("N134").Value = -114121.336

$ws_LTW.Range("H136").Value = 6122.875
$ws_LTW.Range("I136").Value = 8246
$ws_LTW.Range("K136").Value = 24738
$ws_LTW.Range("M136").Value = -22188

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H22").Value = 14999
$ws_WVR.Range("J22").Value = 14999
$ws_WVR.Range("L22").Value = 14999
$ws_WVR.Range("N22").Value = -15585
